# Apply the two changes described by the commit:
#   1. Slide 5's table switches to a different built-in table style
#      (tableStyleId {38D731BC-...} -> {956C350A-...}).
#   2. The presentation's theme colour scheme is swapped from the
#      custom "Integral / Red Violet" palette to the stock
#      "Office Theme / Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$s = $p.Slides.Item(5)
$tbl = $s.Shapes.Item(2).Table
$tbl.ApplyStyle("{956C350A-FF95-4350-99E1-704089984A16}")

# --- 2. Theme colour scheme ------------------------------------------------
# ThemeColorScheme.Item(index).RGB uses OLE COLORREF (BGR) ordering, and the
# index order matches: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink,
# 12 folHlink. Replace the current "Red Violet" colours with the default
# "Office" colours.
$cs = $p.SlideMaster.Theme.ThemeColorScheme
$cs.Item(1).RGB  = 0        # dk1      000000
$cs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388  # dk2      44546A
$cs.Item(4).RGB  = 15132391 # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939 # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501  # accent2  ED7D31
$cs.Item(7).RGB  = 10855845 # accent3  A5A5A5
$cs.Item(8).RGB  = 49407    # accent4  FFC000
$cs.Item(9).RGB  = 12874308 # accent5  4472C4
$cs.Item(10).RGB = 4697456  # accent6  70AD47
$cs.Item(11).RGB = 12673797 # hlink    0563C1
$cs.Item(12).RGB = 7491477  # folHlink 954F72
